$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Change_Risk_Segment" column (F) was removed from the dataset; deleting
# it shifts every later column one slot to the left (G->F, H->G, ... AN->AM).
$ws.Columns("F").Delete() | Out-Null

# Drop the leftover formatted-but-empty tail columns (old AP:AT, now AO:AS)
# so the sheet's used range/column-width metadata matches the new narrower
# table (data now only goes through column AN, which keeps the old AO slot).
$ws.Range("AO1:AS1").EntireColumn.Delete() | Out-Null

# Append the new "Risk_Smoking_Tobacco" header/column at the end.
$ws.Range("AN1").Value = "Risk_Smoking_Tobacco"

# Replace the single data row with the new sample record.
$ws.Range("A2").Value = "West"
$ws.Range("B2").Value = "N"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "N"
$ws.Range("E2").Value = "N"
$ws.Range("F2").Value = "Adherent"
$ws.Range("G2").Value = "N"
$ws.Range("H2").Value = "Y"
$ws.Range("I2").Value = "N"
$ws.Range("J2").Value = "N"
$ws.Range("K2").Value = "N"
$ws.Range("L2").Value = "N"
$ws.Range("M2").Value = "Y"
$ws.Range("N2").Value = "N"
$ws.Range("O2").Value = "Y"
$ws.Range("P2").Value = "Y"
$ws.Range("Q2").Value = "N"
$ws.Range("R2").Value = "N"
$ws.Range("S2").Value = "Y"
$ws.Range("T2").Value = "N"
$ws.Range("U2").Value = "Y"
$ws.Range("V2").Value = "Y"
$ws.Range("W2").Value = "N"
$ws.Range("X2").Value = "N"
$ws.Range("Y2").Value = "N"
$ws.Range("Z2").Value = "N"
$ws.Range("AA2").Value = "N"
$ws.Range("AB2").Value = "N"
$ws.Range("AC2").Value = "N"
$ws.Range("AD2").Value = "N"
$ws.Range("AE2").Value = "N"
$ws.Range("AF2").Value = "N"
$ws.Range("AG2").Value = "N"
$ws.Range("AH2").Value = "N"
$ws.Range("AI2").Value = "N"
$ws.Range("AJ2").Value = "N"
$ws.Range("AK2").Value = "N"
$ws.Range("AL2").Value = "N"
$ws.Range("AM2").Value = 5
$ws.Range("AN2").Value = "N"

# Match the saved selection state.
$ws.Range("B2").Select() | Out-Null
